$d = $word.ActiveDocument

# The "Version History" paragraph currently reads:
#   v1.1 - fixed CSV output to UTF-8, , fixed E.164 mask handling.
# i.e. it has a duplicated ", " before "fixed E.164 mask handling" (the
# runs are: "v1.1 ... UTF-8" | ", " | ", fixed E.164 mask handling" | "."
# followed by the (hidden) "_GoBack" bookmark). The fix removes the extra
# leading ", " from the third run and relocates the "_GoBack" bookmark so
# it now sits right after the (kept) ", " run instead of at the very end
# of the paragraph.

$findRange = $d.Content
$findRange.Find.Execute(", , fixed E.164 mask handling.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $findRange.Find.Found) {
    Write-Output "ERROR: could not locate the duplicated-comma text"
} else {
    $blockStart = $findRange.Start
    $blockEnd = $findRange.End

    # First ", " (2 chars) is kept; the duplicate ", " right after it (also
    # 2 chars) is the part that needs to be removed.
    $keepEnd = $blockStart + 2
    $dupEnd = $keepEnd + 2
    # The trailing "." is the last character of the matched block.
    $periodStart = $blockEnd - 1

    # Relocate the "_GoBack" bookmark to just after the kept ", " run.
    $goBack = $d.Bookmarks.Item("_GoBack")
    $goBack.Delete()
    $d.Bookmarks.Add("_GoBack", $d.Range($keepEnd, $keepEnd))

    # Drop a throwaway bookmark right before the trailing "." run. This
    # keeps that run from being coalesced into the run we are about to
    # edit, mirroring the original document's run layout (the "." stays
    # its own run, separate from "fixed E.164 mask handling").
    $barrierName = "ZZZTempBarrier"
    $d.Bookmarks.Add($barrierName, $d.Range($periodStart, $periodStart))

    # Remove the duplicated ", " text.
    $dupRange = $d.Range($keepEnd, $dupEnd)
    $dupRange.Delete()

    # Remove the throwaway barrier bookmark now that the edit is done.
    $d.Bookmarks.Item($barrierName).Delete()

    Write-Output "Fixed paragraph text: [$($d.Paragraphs.Item(13).Range.Text)]"
}
